# Natmi following Dr Hou advice
# Rewrite the LR-pair data rows (sending/target cluster now includes "ECs")
# and refresh every numeric NATMI statistic (columns E-T) for rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Plxnb2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema4g"
$ws.Range("C2").Value = "Plxnb2"
$ws.Range("D2").Value = "Plxnb2"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1033863333333333
$ws.Range("H2").Value = 0.310159
$ws.Range("I2").Value = 0.05941750124617003
$ws.Range("J2").Value = 0.05941750124617002
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.828998666666665
$ws.Range("N2").Value = 26.486996
$ws.Range("O2").Value = 0.1794455804823882
$ws.Range("P2").Value = 0.1794455804823882
$ws.Range("Q2").Value = 0.9127977991515555
$ws.Range("R2").Value = 8.215180192364
$ws.Range("S2").Value = 0.010662208001932
$ws.Range("T2").Value = 0.010662208001932

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema4g"
$ws.Range("C3").Value = "Plxnb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1033863333333333
$ws.Range("H3").Value = 0.310159
$ws.Range("I3").Value = 0.05941750124617003
$ws.Range("J3").Value = 0.05941750124617002
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 19.33828433333333
$ws.Range("N3").Value = 58.014853
$ws.Range("O3").Value = 0.3930422677296217
$ws.Range("P3").Value = 0.3930422677296217
$ws.Range("Q3").Value = 1.999314310180778
$ws.Range("R3").Value = 17.993828791627
$ws.Range("S3").Value = 0.02335358943262229
$ws.Range("T3").Value = 0.02335358943262229

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema4g"
$ws.Range("C4").Value = "Plxnb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1033863333333333
$ws.Range("H4").Value = 0.310159
$ws.Range("I4").Value = 0.05941750124617003
$ws.Range("J4").Value = 0.05941750124617002
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 21.03425566666667
$ws.Range("N4").Value = 63.102767
$ws.Range("O4").Value = 0.4275121517879902
$ws.Range("P4").Value = 0.4275121517879902
$ws.Range("Q4").Value = 2.174654567772556
$ws.Range("R4").Value = 19.571891109953
$ws.Range("S4").Value = 0.02540170381161574
$ws.Range("T4").Value = 0.02540170381161573

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema4g"
$ws.Range("C5").Value = "Plxnb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.264650666666667
$ws.Range("H5").Value = 3.793952
$ws.Range("I5").Value = 0.7268115633849387
$ws.Range("J5").Value = 0.7268115633849388
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.828998666666665
$ws.Range("N5").Value = 26.486996
$ws.Range("O5").Value = 0.1794455804823882
$ws.Range("P5").Value = 0.1794455804823882
$ws.Range("Q5").Value = 11.16559904979911
$ws.Range("R5").Value = 100.490391448192
$ws.Range("S5").Value = 0.1304231228929224
$ws.Range("T5").Value = 0.1304231228929224

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema4g"
$ws.Range("C6").Value = "Plxnb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.264650666666667
$ws.Range("H6").Value = 3.793952
$ws.Range("I6").Value = 0.7268115633849387
$ws.Range("J6").Value = 0.7268115633849388
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.33828433333333
$ws.Range("N6").Value = 58.014853
$ws.Range("O6").Value = 0.3930422677296217
$ws.Range("P6").Value = 0.3930422677296217
$ws.Range("Q6").Value = 24.45617417433956
$ws.Range("R6").Value = 220.105567569056
$ws.Range("S6").Value = 0.285667665084928
$ws.Range("T6").Value = 0.285667665084928

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema4g"
$ws.Range("C7").Value = "Plxnb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.264650666666667
$ws.Range("H7").Value = 3.793952
$ws.Range("I7").Value = 0.7268115633849387
$ws.Range("J7").Value = 0.7268115633849388
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 21.03425566666667
$ws.Range("N7").Value = 63.102767
$ws.Range("O7").Value = 0.4275121517879902
$ws.Range("P7").Value = 0.4275121517879902
$ws.Range("Q7").Value = 26.60098545168711
$ws.Range("R7").Value = 239.408869065184
$ws.Range("S7").Value = 0.3107207754070884
$ws.Range("T7").Value = 0.3107207754070884

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sema4g"
$ws.Range("C8").Value = "Plxnb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.371961
$ws.Range("H8").Value = 1.115883
$ws.Range("I8").Value = 0.2137709353688912
$ws.Range("J8").Value = 0.2137709353688912
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.828998666666665
$ws.Range("N8").Value = 26.486996
$ws.Range("O8").Value = 0.1794455804823882
$ws.Range("P8").Value = 0.1794455804823882
$ws.Range("Q8").Value = 3.284043173051999
$ws.Range("R8").Value = 29.556388557468
$ws.Range("S8").Value = 0.03836024958753376
$ws.Range("T8").Value = 0.03836024958753376

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sema4g"
$ws.Range("C9").Value = "Plxnb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.371961
$ws.Range("H9").Value = 1.115883
$ws.Range("I9").Value = 0.2137709353688912
$ws.Range("J9").Value = 0.2137709353688912
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 19.33828433333333
$ws.Range("N9").Value = 58.014853
$ws.Range("O9").Value = 0.3930422677296217
$ws.Range("P9").Value = 0.3930422677296217
$ws.Range("Q9").Value = 7.193087578911
$ws.Range("R9").Value = 64.737788210199
$ws.Range("S9").Value = 0.0840210132120714
$ws.Range("T9").Value = 0.0840210132120714

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema4g"
$ws.Range("C10").Value = "Plxnb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.371961
$ws.Range("H10").Value = 1.115883
$ws.Range("I10").Value = 0.2137709353688912
$ws.Range("J10").Value = 0.2137709353688912
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 21.03425566666667
$ws.Range("N10").Value = 63.102767
$ws.Range("O10").Value = 0.4275121517879902
$ws.Range("P10").Value = 0.4275121517879902
$ws.Range("Q10").Value = 7.823922772028999
$ws.Range("R10").Value = 70.415304948261
$ws.Range("S10").Value = 0.09138967256928607
$ws.Range("T10").Value = 0.09138967256928607

